$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The update swaps the full data (columns B:AC) between each of these
# row pairs, while the row index in column A stays fixed per row.
$pairs = @(
    @(135, 136),
    @(139, 140),
    @(142, 143),
    @(144, 145),
    @(214, 215)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AC$r1")
    $rng2 = $ws.Range("B$r2`:AC$r2")

    $val1 = $rng1.Value()
    $val2 = $rng2.Value()

    $rng1.Value = $val2
    $rng2.Value = $val1
}
